$d = $word.ActiveDocument

# 1. Title (appears twice: Heading1 at top, and bold run near bottom) - replace all occurrences
$d.Content.Find.Execute(
    "Play Miss Kitty Free: Review and Strategies | Top Online Slots", $true, $false, $false, $false, $false,
    $true, 1, $false, "Play Miss Kitty Free | Review and Gameplay", 2
)

# 2. "What we like" bullet list items
$d.Content.Find.Execute(
    "Sticky Wilds feature during Free Spins rounds", $true, $false, $false, $false, $false,
    $true, 1, $false, "Simple yet legible graphics", 2
)

$d.Content.Find.Execute(
    "Up to 50 pay lines available", $true, $false, $false, $false, $false,
    $true, 1, $false, "Multiple bonus features", 2
)

$d.Content.Find.Execute(
    "Simple graphics and legible symbols", $true, $false, $false, $false, $false,
    $true, 1, $false, "Sticky Wilds increase winning chances", 2
)

$d.Content.Find.Execute(
    "Gamble feature to increase winnings", $true, $false, $false, $false, $false,
    $true, 1, $false, "Vibrant and colorful design", 2
)

# 3. "What we don't like" bullet list items
$d.Content.Find.Execute(
    "Lower payouts compared to other popular slot games", $true, $false, $false, $false, $false,
    $true, 1, $false, "Limited variety of symbols", 2
)

$d.Content.Find.Execute(
    "No progressive jackpot feature", $true, $false, $false, $false, $false,
    $true, 1, $false, "Gamble feature may not appeal to all players", 2
)

# 4. Final italic summary paragraph
$d.Content.Find.Execute(
    "Play Miss Kitty for free and read our review of the game, including winning strategies, bonus features, and winning opportunities. Top online slot gaming!",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "Discover the charm of Miss Kitty slot game and play for free. Read our review to learn more.", 2
)
